# VO_AU_CONFIG.xlsx edit
#
# Commit: "re-write json config, copy xlsx config, re-structure DatabaseExtract.py"
#
# The xlsx-relevant part of that commit:
#   - duplicate the CARSGUIDE worksheet into a new worksheet named
#     "YURY.CARSGUIDE" placed right after AUTOTRADER (last tab);
#   - the VO_AU sheet had row 47 grow a touch taller (13.8 -> 15) and the
#     window scrolled/selected further down (selection moved from C42 to
#     C49, i.e. roughly the same relative offset it had before, now that
#     more rows are in view);
#   - the AUTOTRADER sheet's remembered selection was reset back to A1.

$wb = $excel.ActiveWorkbook

# --- 1. Copy CARSGUIDE -> new "YURY.CARSGUIDE" tab, placed after AUTOTRADER ---
$carsguide = $wb.Worksheets.Item("CARSGUIDE")
$autotrader = $wb.Worksheets.Item("AUTOTRADER")

# Copy() with an "after" target sheet inserts the duplicate right after it,
# mirroring how the new tab ends up last (after AUTOTRADER) in the workbook.
$carsguide.Copy($null, $autotrader)

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "YURY.CARSGUIDE"

# --- 2. VO_AU sheet tweaks ---
$voau = $wb.Worksheets.Item("VO_AU")
$voau.Activate()

# Row 47 (expect_column_values_to_be_in_set_2 / TYPE) got a touch taller.
$voau.Rows.Item(47).RowHeight = 15

# Selection/scroll moved on from C42 to C49.
$voau.Range("C49").Select() | Out-Null

# --- 3. AUTOTRADER sheet: selection reset to A1 ---
$autotrader.Activate()
$autotrader.Range("A1").Select() | Out-Null

# Leave the newly-added sheet as the active tab afterwards, matching a
# "just finished copying it" end state.
$newSheet.Activate()
